# WSN_2009.xlsx: add season-record columns (Wins / Losses / Ties)
# after the existing "Unnamed: 28" column (AC), filling every data row
# (2-57) with the team's 2009 record: 59 wins, 103 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - clone the formatting of the last existing
# header cell (AC1, style index 1: bold font + thin border + centered/top
# alignment) so the new headers look identical to the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player/data row.
$firstRow = 2
$lastRow = 57
$wins = 59
$losses = 103
$ties = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
